$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row for new columns D:I
$ws.Range("D1").Value = "NegCueSD"
$ws.Range("E1").Value = "NeuCueSD"
$ws.Range("F1").Value = "AmbCueSD"
$ws.Range("G1").Value = "NegCueSR"
$ws.Range("H1").Value = "NeuCueSR"
$ws.Range("I1").Value = "AmbCueSR"

# Data rows 2..19 for columns D:I
$data = @(
    @(0.072172, -1.04497, -0.636703, -0.836327, -0.158059, -0.0311363),
    @(-0.353648, -1.35444, -0.51698, -0.68824, -0.784818, 0.0564267),
    @(-0.256401, -0.856722, 0.102033, -0.401194, -0.0213203, -0.335401),
    @(0.0495691, -0.417189, 0.117841, -0.829303, -0.7998, -1.25495),
    @(0.165005, 1.75244, 0.89788, 0.0306742, -0.506285, 0.30956),
    @(-0.00315595, 0.511022, -0.122243, -1.58629, -1.40585, -1.14161),
    @(0.601673, 0.0996013, -0.339815, -1.46858, -1.05149, -1.09291),
    @(-1.53888, -1.32329, -1.2447, -2.91198, -2.19771, -0.574445),
    @(-0.639092, -1.9363, -1.68345, -0.598612, -0.911635, -1.09657),
    @(0.305705, 0.444902, -0.531341, -0.86239, -1.20402, -0.173988),
    @(0.558604, -0.138457, -0.0781469, -0.274484, -0.758659, -0.173326),
    @(-1.55924, -1.77223, -1.41806, -1.60245, -1.692, -0.834505),
    @(-0.207432, 0.443186, 0.166934, 0.787154, -0.361806, -0.266177),
    @(0.406919, -0.645518, -0.905686, -1.10684, -1.30483, -0.875592),
    @(-0.351634, 0.165144, -0.219244, -0.329889, -0.932411, -0.469584),
    @(0.993676, 0.335478, 0.414095, -0.593378, -0.523391, -0.45518),
    @(0.436064, -0.102812, -0.413867, -1.12862, -0.933742, -0.225554),
    @(-0.0722588, 0.357688, -0.906409, -0.866392, -0.617351, -1.37932)
)

$r = 2
foreach ($row in $data) {
    $ws.Range("D$r").Value = $row[0]
    $ws.Range("E$r").Value = $row[1]
    $ws.Range("F$r").Value = $row[2]
    $ws.Range("G$r").Value = $row[3]
    $ws.Range("H$r").Value = $row[4]
    $ws.Range("I$r").Value = $row[5]
    $r++
}

# Update the selection to match the diff (D1:I19, active cell D1)
$ws.Range("D1:I19").Select()
